$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.2030812324929972
    "C2" = 0.5098039215686274
    "J2" = 0.02240896358543417
    "O2" = 0.001400560224089636
    "P2" = 0.1666666666666667
    "S2" = 0.09663865546218488
    "B3" = 0.01401869158878505
    "C3" = 0.03271028037383177
    "J3" = 0.0514018691588785
    "P3" = 0.7336448598130841
    "S3" = 0.1682242990654206
    "J4" = 0.1071428571428571
    "O4" = 0.008928571428571428
    "P4" = 0.5982142857142857
    "S4" = 0.2857142857142857
    "P5" = 0.75
    "S5" = 0.25
    "B6" = 0.072265625
    "D6" = 0.01171875
    "E6" = 0.00390625
    "F6" = 0.060546875
    "J6" = 0.267578125
    "O6" = 0.025390625
    "Q6" = 0.185546875
    "R6" = 0.044921875
    "S6" = 0.328125
    "B7" = 0.1350114416475973
    "D7" = 0.009153318077803204
    "F7" = 0.07322654462242563
    "J7" = 0.1327231121281464
    "O7" = 0.02745995423340961
    "Q7" = 0.1624713958810069
    "R7" = 0.07780320366132723
    "S7" = 0.3821510297482837
    "B8" = 0.102803738317757
    "D8" = 0.02388369678089304
    "E8" = 0.002076843198338525
    "F8" = 0.06542056074766354
    "J8" = 0.1121495327102804
    "O8" = 0.02596053997923157
    "Q8" = 0.1962616822429906
    "R8" = 0.0778816199376947
    "S8" = 0.3935617860851506
    "B9" = 0.08484848484848485
    "D9" = 0.02121212121212121
    "F9" = 0.06060606060606061
    "J9" = 0.1272727272727273
    "O9" = 0.03636363636363636
    "Q9" = 0.1757575757575758
    "R9" = 0.09090909090909091
    "S9" = 0.403030303030303
    "B10" = 0.1070567986230637
    "D10" = 0.023407917383821
    "E10" = 0.001721170395869191
    "F10" = 0.06506024096385542
    "J10" = 0.13184165232358
    "O10" = 0.03132530120481928
    "Q10" = 0.2351118760757315
    "R10" = 0.07091222030981068
    "S10" = 0.3335628227194492
    "G11" = 0.1376281112737921
    "J11" = 0.102489019033675
    "K11" = 0.1976573938506589
    "L11" = 0.5490483162518301
    "S11" = 0.01317715959004392
    "G12" = 0.7131979695431472
    "J12" = 0.2106598984771574
    "K12" = 0.01015228426395939
    "L12" = 0.03045685279187817
    "S12" = 0.03553299492385787
    "G13" = 0.7634408602150538
    "J13" = 0.1397849462365591
    "S13" = 0.09677419354838709
    "F15" = 0.03119266055045872
    "H15" = 0.1541284403669725
    "I15" = 0.04220183486238532
    "J15" = 0.3486238532110092
    "K15" = 0.07522935779816514
    "M15" = 0.009174311926605505
    "N15" = 0.003669724770642202
    "O15" = 0.04954128440366973
    "S15" = 0.2862385321100918
    "F16" = 0.02100840336134454
    "H16" = 0.1680672268907563
    "I16" = 0.04831932773109244
    "J16" = 0.3970588235294117
    "K16" = 0.1134453781512605
    "M16" = 0.01260504201680672
    "O16" = 0.05042016806722689
    "S16" = 0.1890756302521008
    "F17" = 0.02228163992869875
    "H17" = 0.1809269162210339
    "I17" = 0.0659536541889483
    "J17" = 0.4322638146167558
    "K17" = 0.09982174688057041
    "M17" = 0.01693404634581105
    "N17" = 0.00089126559714795
    "O17" = 0.06417112299465241
    "S17" = 0.1167557932263815
    "F18" = 0.01570680628272251
    "H18" = 0.143979057591623
    "I18" = 0.08900523560209424
    "J18" = 0.4371727748691099
    "K18" = 0.0968586387434555
    "M18" = 0.005235602094240838
    "N18" = 0.002617801047120419
    "O18" = 0.08638743455497382
    "S18" = 0.1230366492146597
    "F19" = 0.0233433734939759
    "H19" = 0.2085843373493976
    "J19" = 0.376882530120482
    "K19" = 0.1204819277108434
    "M19" = 0.02371987951807229
    "N19" = 0.0007530120481927711
    "O19" = 0.07341867469879518
    "S19" = 0.1065512048192771
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

Write-Host "Applied $($values.Count) cell updates"